$wb = $excel.ActiveWorkbook

# --- Add the new "pageRank" worksheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "pageRank"

# Column A width to match source formatting (closest value reachable
# through the host's pixel-quantized ColumnWidth setter; the host rounds
# to the nearest 1/6-character increment, so 36.21875 itself is not
# reachable bit-for-bit — 35.25 lands on the nearest supported width).
$ws.Columns.Item(1).ColumnWidth = 35.25

# --- Header rows ---
$ws.Range("A1").Value = "total vertices is: 172"
$ws.Range("A2").Value = "total edges is: 370"

# --- PageRank values (row 4 onward; row 3 intentionally left blank) ---
$pageRankValues = @(
  "CAYWK: 0.6085755960976994"
  "CAYWJ: 0.24109125671940038"
  "CAYWG: 0.9800828853228987"
  "CAYFS: 0.1513954755488643"
  "CAYFO: 0.2395943184668877"
  "CAYGK: 0.1348166008352925"
  "CAYOP: 0.20980416011955139"
  "CAYGH: 0.25304459749626695"
  "CAYOJ: 0.20980416898862125"
  "CAYWP: 0.30664273353511834"
  "CAYWL: 0.2215440640249542"
  "CAYXL: 0.879826234881961"
  "CAYGZ: 0.33333336828892207"
  "CAYPC: 0.18115288701701868"
  "CAYXJ: 0.32443794629794653"
  "CAYGX: 0.1836889702227043"
  "CAYGW: 0.4327445857283516"
  "CAYXH: 0.13998273091512198"
  "CAYGV: 0.221346822657079"
  "CAYGT: 0.33886320849820173"
  "CAYXE: 0.29128398340571926"
  "CAYXC: 0.2218136070089129"
  "CAYGR: 0.2511564793942326"
  "CAYGP: 0.19452655083024212"
  "CAYOW: 0.7105370542429214"
  "CAYGL: 0.24940199844009867"
  "CAYHK: 0.277092121785632"
  "CAYPR: 0.14400363385408516"
  "CAYXY: 0.3095818510469943"
  "CAYHI: 0.34046636521171436"
  "CAYXX: 0.17635543188865593"
  "CAYPO: 0.23326652991324054"
  "CAYPN: 0.22134681682250662"
  "CAYXU: 0.24406604444276989"
  "CAYPM: 0.3143536336219239"
  "CAYXT: 0.2443317318632756"
  "CAYHD: 0.3000195949749606"
  "CAYPL: 0.2146592467318236"
  "CAYXS: 0.4950554653262319"
  "CAZRJ: 0.23280492060598407"
  "CAYPH: 0.27124507758043803"
  "CAYXP: 0.30590729149023016"
  "CAYXN: 0.2626484599232371"
  "CAYUB: 0.18115288701701868"
  "CAYDQ: 0.20476516465056047"
  "CAYDP: 0.1474129387184554"
  "CAYLW: 0.4863502573676459"
  "CAZFN: 0.256273660487246"
  "CAZFM: 0.18115288701701868"
  "CAYTZ: 0.43109186422645307"
  "CAYUY: 0.2390390728280849"
  "CAYUX: 0.24002808095617234"
  "CAYEG: 1.1431420305967814"
  "CAYMO: 0.32154885154415197"
  "CAYMM: 0.376522265245143"
  "CAYUT: 0.26314611517523223"
  "CAYUL: 1.2690404400947906"
  "CAYNA: 0.24119960707160387"
  "CAZPB: 0.15132319703478106"
  "CAYEV: 0.8115288701701867"
  "CAYER: 0.2837310000870396"
  "CAYVB: 0.33677725958358407"
  "CAYMT: 0.2066416958928582"
  "CAYEK: 0.2569086486499763"
  "CAYFJ: 0.1513954755488643"
  "CAYVZ: 0.2578721481347357"
  "CAYFH: 0.23151200016385393"
  "CAYFC: 0.23796120802036166"
  "CAYVR: 1.9487323563951993"
  "CAYFB: 0.783771415370468"
  "CAYVQ: 0.41080977416713327"
  "CAYFA: 0.29250298169761113"
  "CAYVP: 0.6401097420062235"
  "CAYVO: 0.3107755334439678"
  "CAYVM: 0.20706755202158056"
  "CAYNC: 0.21240573449222372"
  "CAYSG: 0.1513954755488643"
  "CAYBR: 0.13998273091512198"
  "CAYSB: 0.36992592686258585"
  "CAYBL: 0.19447827622567634"
  "CAYJT: 0.18452144768075524"
  "CAYBK: 0.32918239628103474"
  "CAYBG: 0.30105084522415415"
  "CAYCG: 0.1839863647692071"
  "CAYKL: 0.2745359735825577"
  "CAYCD: 0.1839863647692071"
  "CAYCB: 0.34881037605590104"
  "CAYKG: 0.15600961944172737"
  "CAYKF: 0.13998273091512198"
  "CAYSM: 0.31401067801155136"
  "CAZUM: 0.18614047665194536"
  "CAYSK: 0.1865951989504465"
  "CAYSJ: 0.20480281215569202"
  "CAYKA: 0.22730372932996404"
  "CAYBX: 0.3010499025332129"
  "CAYTH: 0.40783417195730176"
  "CAYTE: 0.22140304072034572"
  "CAYCS: 0.26314613345597065"
  "CAYCO: 0.29187945365274837"
  "CAYCL: 0.174547275866722"
  "CAZEM: 0.32115924140635344"
  "CAZMT: 0.14400363385408516"
  "CAYSY: 0.260595038899752"
  "CAYKQ: 0.4141441286039471"
  "CAYDF: 0.3220642776296798"
  "CAYLL: 0.13998273091512198"
  "CAYTS: 0.4002649327605129"
  "CAYTQ: 0.15600960242554457"
  "CAYLH: 0.2525791809442863"
  "CAYTL: 0.28373100413219743"
  "CAYLC: 0.16095999897325863"
  "CAYCY: 0.25485339806492513"
  "CAYQD: 0.2395943111983114"
  "CAYHZ: 0.7675658128549724"
  "CAYQB: 0.7184672215063099"
  "CAYHY: 0.24272005691868254"
  "CAYYJ: 0.2930029233420861"
  "CAYYH: 0.27709213040136443"
  "CAYYG: 0.20480281215569202"
  "CAYYF: 0.14400363385408516"
  "CAYHU: 0.27633162665192035"
  "CAYYE: 0.2604064971190653"
  "CAYYD: 0.24433172450489474"
  "CAYYC: 1.8277819846912897"
  "CAYHR: 0.27806906687794997"
  "CAYYB: 0.17797462563592753"
  "CAYPY: 0.19887946626395483"
  "CAYPX: 0.2774939563017271"
  "CAYPW: 0.14400363385408516"
  "CAZBF: 0.1813152989154998"
  "CAYHM: 0.25165529514783724"
  "CAYQU: 0.17635543188865593"
  "CAYQT: 0.8006008320816873"
  "CAYIK: 0.33333336828892207"
  "CAYAC: 0.23280492448915147"
  "CAYYZ: 1.442402034604976"
  "CAYQR: 0.29128398340571926"
  "CAYQQ: 0.28842646242172776"
  "CAYYY: 0.34339566367742336"
  "CAYIF: 0.2675691554368246"
  "CAYQM: 0.2820008951853523"
  "CAYYU: 0.18980685106235584"
  "CAYYT: 0.5069886275905724"
  "CAYQL: 0.13998273091512198"
  "CAZKE: 0.33564258761597165"
  "CAYQK: 0.24745170879134812"
  "CAYYR: 0.6095949263801407"
  "CAYYQ: 0.2887305034297961"
  "CAZSJ: 0.314353623468964"
  "CAYQG: 0.21251988038596165"
  "CAYQF: 0.2217375707667611"
  "CAYAT: 0.3807615140378301"
  "CAYRB: 0.33333335780224543"
  "CAYRA: 0.1513954755488643"
  "CAYZG: 0.33333335780224543"
  "CAYZF: 1.2481758347581324"
  "CAYAM: 0.2588044911214399"
  "CAYQZ: 0.2215440562628191"
  "CAYQY: 0.17063633876852458"
  "CAYQX: 0.22266512546428752"
  "CAYIO: 0.2682667809531727"
  "CAYAG: 0.25217122405356807"
  "CAYRT: 0.7086538498122924"
  "CAYBC: 0.24888021063193808"
  "CAYZV: 0.6267918413218422"
  "CAYZT: 0.14400363385408516"
  "CAYRL: 0.4101289084477824"
  "CAYZS: 0.3453887430116991"
  "CAYZR: 0.1348166008352925"
  "CAYAY: 0.2570903721817179"
  "CAYAX: 0.28373100894891057"
  "CAYZP: 0.14400363385408516"
)

$arr = New-Object 'object[,]' $pageRankValues.Length,1
for ($i = 0; $i -lt $pageRankValues.Length; $i++) {
    $arr[$i,0] = $pageRankValues[$i]
}
$startRow = 4
$endRow = $startRow + $pageRankValues.Length - 1
$ws.Range("A$startRow`:A$endRow").Value = $arr

# --- Selections ---
$ws.Range("B8").Select()

$clustering = $wb.Worksheets.Item("clustering coefficient")
$clustering.Activate()
$clustering.Range("J32").Select()

$ws.Activate()
